$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-19 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-20 Friday", 2) | Out-Null
$d.Content.Find.Execute("861÷9=95, 6", $true, $false, $false, $false, $false, $true, 1, $false, "831÷9=92, 3", 2) | Out-Null
$d.Content.Find.Execute("639÷5=127, 4", $true, $false, $false, $false, $false, $true, 1, $false, "869÷9=96, 5", 2) | Out-Null
$d.Content.Find.Execute("871÷8=108, 7", $true, $false, $false, $false, $false, $true, 1, $false, "775÷9=86, 1", 2) | Out-Null
$d.Content.Find.Execute("321÷9=35, 6", $true, $false, $false, $false, $false, $true, 1, $false, "824÷3=274, 2", 2) | Out-Null
$d.Content.Find.Execute("218÷9=24, 2", $true, $false, $false, $false, $false, $true, 1, $false, "494÷6=82, 2", 2) | Out-Null
$d.Content.Find.Execute("328÷3=109, 1", $true, $false, $false, $false, $false, $true, 1, $false, "532÷7=76, 0", 2) | Out-Null
$d.Content.Find.Execute("749÷4=187, 1", $true, $false, $false, $false, $false, $true, 1, $false, "627÷7=89, 4", 2) | Out-Null
$d.Content.Find.Execute("433÷5=86, 3", $true, $false, $false, $false, $false, $true, 1, $false, "501÷5=100, 1", 2) | Out-Null
$d.Content.Find.Execute("192÷3=64, 0", $true, $false, $false, $false, $false, $true, 1, $false, "942÷8=117, 6", 2) | Out-Null
$d.Content.Find.Execute("100÷6=16, 4", $true, $false, $false, $false, $false, $true, 1, $false, "969÷3=323, 0", 2) | Out-Null
$d.Content.Find.Execute("145÷7=20, 5", $true, $false, $false, $false, $false, $true, 1, $false, "901÷7=128, 5", 2) | Out-Null
$d.Content.Find.Execute("102÷2=51, 0", $true, $false, $false, $false, $false, $true, 1, $false, "576÷4=144, 0", 2) | Out-Null
$d.Content.Find.Execute("268÷3=89, 1", $true, $false, $false, $false, $false, $true, 1, $false, "151÷7=21, 4", 2) | Out-Null
$d.Content.Find.Execute("511÷6=85, 1", $true, $false, $false, $false, $false, $true, 1, $false, "674÷4=168, 2", 2) | Out-Null
$d.Content.Find.Execute("810÷7=115, 5", $true, $false, $false, $false, $false, $true, 1, $false, "985÷3=328, 1", 2) | Out-Null
$d.Content.Find.Execute("645÷8=80, 5", $true, $false, $false, $false, $false, $true, 1, $false, "805÷8=100, 5", 2) | Out-Null
$d.Content.Find.Execute("478÷7=68, 2", $true, $false, $false, $false, $false, $true, 1, $false, "487÷8=60, 7", 2) | Out-Null
$d.Content.Find.Execute("245÷9=27, 2", $true, $false, $false, $false, $false, $true, 1, $false, "651÷8=81, 3", 2) | Out-Null
$d.Content.Find.Execute("166÷9=18, 4", $true, $false, $false, $false, $false, $true, 1, $false, "557÷8=69, 5", 2) | Out-Null
$d.Content.Find.Execute("526÷6=87, 4", $true, $false, $false, $false, $false, $true, 1, $false, "990÷6=165, 0", 2) | Out-Null
$d.Content.Find.Execute("600÷7=85, 5", $true, $false, $false, $false, $false, $true, 1, $false, "263÷8=32, 7", 2) | Out-Null
$d.Content.Find.Execute("962÷3=320, 2", $true, $false, $false, $false, $false, $true, 1, $false, "691÷8=86, 3", 2) | Out-Null
$d.Content.Find.Execute("527÷7=75, 2", $true, $false, $false, $false, $false, $true, 1, $false, "950÷9=105, 5", 2) | Out-Null
$d.Content.Find.Execute("162÷8=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "378÷4=94, 2", 2) | Out-Null
$d.Content.Find.Execute("789÷7=112, 5", $true, $false, $false, $false, $false, $true, 1, $false, "634÷2=317, 0", 2) | Out-Null

$d.Save()
